$S = @(
    "Ementa atual:",
    "Ementa modificada (dados modificados em vermelho):",
    "LOB1053",
    "Nome:",
    " Física III",
    "Name:",
    "Physics III",
    "Créditos-aula:",
    "4",
    "Créditos-trabalho",
    "0",
    "Carga horária:",
    "60 h",
    "Ativação:",
    "01/01/2018",
    "Semestre ideal:",
    "EF-3,EM-3,EA-4,EB-4,EP-4,EQD-3,EQN-4",
    "Objetivos:",
    "2342277 - Bertha María Cuadros Melgar",
    "Objectives:",
    "To introduce to students the basic concepts of electromagnetism such as electric charge, electric field, electric potential, magnetic field, and Lorentz force showing their applications to several devices and configurations. In addition, the students are going to get familiarized with Gauss, Ampère, and Faraday laws. Finally, students should understand the relation between magnetic and electric fields and how to generate electric current from a magnetic field through induction.",
    "Docentes responsáveis:",
    "Programa resumido:",
    "Semestral",
    "Short syllabus:",
    "Electric Charge and Matter. Electric fields. The Gauss' law . Electric Potential . Capacitors and Dielectrics. Electric Current and Resistance. Magnetic Fields . Magnetic Fields sources. Electromagnetic induction and inductance . Magnetic Properties of Matter. Maxwell's equations.",
    "Programa:",
    "Syllabus:",
    "1) Electric charge and electric force: electric charge; conductors and insulators; Coulomb's law; quantization and conservation.2) Electric field: concepts; field lines; point charge and dipole, continuous distribution.3) Gauss' law: flow; applications in cylindrical, flat and spherical geometries.4) Electric potential: concept and calculation; energy, potential and electric field, equipotential surfaces; punctual loads, electric dipole and continuous distributions.5) Capacitors and dielectrics: capacitance, energy and calculation, associations, dielectrics.6) Electric current and resistance: current density, resistance and resistivity as a function of temperature; Ohm's law, power, semiconductors and superconductors.7) Magnetic field: Biot-Savart law.8) Ampère's law and applications: magnetic field of a coil, solenoid, and toroids.9) Electromagnetic induction and inductance: Faraday's law, Lenz's law.10) Magnetic properties of matter.11) Maxwell's equations.",
    "Avaliação:",
    "Método:",
    "Critério:",
    "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.",
    "Norma de recuperação:",
    "NF≥ 5,0.",
    "Bibliografia:",
    "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.",
    "Requisitos:",
    "LOB1004 -  Cálculo II  (Requisito fraco)`n",
    "LOB1019 -  Física II  (Requisito fraco)`n"
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlShiftUp = -4162

# ---------------------------------------------------------------
# Row 13: originally only had B13/C13 ("2342277 - Bertha..." = old
# index 22). Target row 13 needs A13 (style like A14/A3, bold) plus
# new values in B13/C13, and a 60pt custom row height.
# ---------------------------------------------------------------
$ws.Range("A3").Copy()
$ws.Range("A13").PasteSpecial($xlPasteFormats)
$ws.Range("A13").Value = $S[22]
$ws.Range("B13").Value = $S[23]
$ws.Range("C13").Value = $S[23]
$ws.Rows.Item(13).RowHeight = 60

# ---------------------------------------------------------------
# Row 14: all three cells already exist - just replace values.
# Height stays 60 (unchanged).
# ---------------------------------------------------------------
$ws.Range("A14").Value = $S[24]
$ws.Range("B14").Value = $S[25]
$ws.Range("C14").Value = $S[25]
$ws.Rows.Item(14).RowHeight = 60

# ---------------------------------------------------------------
# Row 15: all three cells already exist - replace values, height
# changes from 60 to 120.
# ---------------------------------------------------------------
$ws.Range("A15").Value = $S[26]
$ws.Range("B15").Value = $S[14]
$ws.Range("C15").Value = $S[14]
$ws.Rows.Item(15).RowHeight = 120

# ---------------------------------------------------------------
# Row 16: unchanged entirely (left as-is).
# ---------------------------------------------------------------

# ---------------------------------------------------------------
# Row 17: A17 keeps its value ("Avaliação:"); B17/C17 must be
# removed entirely, and the custom row height (120) removed too
# (back to sheet default).
# ---------------------------------------------------------------
$ws.Range("B17:C17").Delete($xlShiftUp)
$ws.Rows.Item(17).AutoFit()

# ---------------------------------------------------------------
# Row 18: originally only had A18. Target needs B18/C18 added
# (style like B19/C19, wrap text), A18 value replaced, and a 60pt
# custom row height.
# ---------------------------------------------------------------
$ws.Range("B19").Copy()
$ws.Range("B18").PasteSpecial($xlPasteFormats)
$ws.Range("C19").Copy()
$ws.Range("C18").PasteSpecial($xlPasteFormats)
$ws.Range("A18").Value = $S[30]
$ws.Range("B18").Value = $S[18]
$ws.Range("C18").Value = $S[18]
$ws.Rows.Item(18).RowHeight = 60

# ---------------------------------------------------------------
# Row 19: all three cells already exist - just replace values.
# Height stays 60.
# ---------------------------------------------------------------
$ws.Range("A19").Value = $S[31]
$ws.Range("B19").Value = $S[32]
$ws.Range("C19").Value = $S[32]
$ws.Rows.Item(19).RowHeight = 60

# ---------------------------------------------------------------
# Row 20: all three cells already exist - just replace values.
# Height stays 60.
# ---------------------------------------------------------------
$ws.Range("A20").Value = $S[33]
$ws.Range("B20").Value = $S[34]
$ws.Range("C20").Value = $S[34]
$ws.Rows.Item(20).RowHeight = 60

# ---------------------------------------------------------------
# Row 21: all three cells already exist - replace values, height
# changes from 60 to 120.
# ---------------------------------------------------------------
$ws.Range("A21").Value = $S[35]
$ws.Range("B21").Value = $S[36]
$ws.Range("C21").Value = $S[36]
$ws.Rows.Item(21).RowHeight = 120

# ---------------------------------------------------------------
# Row 22: A22 value replaced; B22/C22 removed entirely; custom
# row height (120) removed (back to sheet default).
# ---------------------------------------------------------------
$ws.Range("A22").Value = $S[37]
$ws.Range("B22:C22").Delete($xlShiftUp)
$ws.Rows.Item(22).AutoFit()

# ---------------------------------------------------------------
# Row 23: originally only had A23. Target needs A23 removed and
# B23/C23 added (style like B24/C24), with a 30pt custom height.
# ---------------------------------------------------------------
$ws.Range("B24").Copy()
$ws.Range("B23").PasteSpecial($xlPasteFormats)
$ws.Range("C24").Copy()
$ws.Range("C23").PasteSpecial($xlPasteFormats)
$ws.Range("B23").Value = $S[38]
$ws.Range("C23").Value = $S[38]
$ws.Range("A23:A23").Delete($xlShiftUp)
$ws.Rows.Item(23).RowHeight = 30

# ---------------------------------------------------------------
# Row 24: B24/C24 already exist - just replace values. Height
# stays 30.
# ---------------------------------------------------------------
$ws.Range("B24").Value = $S[39]
$ws.Range("C24").Value = $S[39]
$ws.Rows.Item(24).RowHeight = 30

# ---------------------------------------------------------------
# Row 25: no longer exists in the target sheet - remove it
# entirely (shifts nothing since it's the last row).
# ---------------------------------------------------------------
$ws.Rows.Item(25).Delete()

Write-Output "Edit complete"
